# "break out stock.yaml completed"
# 1) D388:D418 were stored as text (inlineStr) BSE codes; re-save them as real
#    numbers (same digits, numeric type).
# 2) Append 12 new "day" rows (419-430) for stocks scraped on 21/08/2024.
#    The dimension (A1:I430) updates automatically as cells are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- 1) Convert D388:D418 from text to numeric, keeping the same digits ---
for ($r = 388; $r -le 418; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $text = $cell.Value2
    $cell.Value = [double]$text
}

# --- 2) Append the 12 new rows ---
$newRows = @(
    @(1,  "SIEMENS",    "Siemens Limited",                                          "500550", -0.01,               7083.05, 148030,  "day", "21/08/2024 11:36:40"),
    @(2,  "M&M",        "Mahindra & Mahindra Limited",                               "500520", -0.07000000000000001, 2769.4, 1497501, "day", "21/08/2024 11:36:40"),
    @(3,  "AXISBANK",   "Axis Bank Limited",                                         "532215", 0.55,                1174.4, 9480630, "day", "21/08/2024 11:36:40"),
    @(4,  "CANFINHOME", "Can Fin Homes Limited",                                     "511196", -0.08,                  847,  352026, "day", "21/08/2024 11:36:40"),
    @(5,  "HINDALCO",   "Hindalco Industries Limited",                               "500440", 1.89,                685.6, 7815238, "day", "21/08/2024 11:36:40"),
    @(6,  "MARICO",     "Marico Limited",                                            "531642", 1.55,                679.3, 1206303, "day", "21/08/2024 11:36:40"),
    @(7,  "GNFC",       "Gujarat Narmada Valley Fertilizers And Chemicals Limited",  "500670", 1.44,               668.15,  940951, "day", "21/08/2024 11:36:40"),
    @(8,  "INDHOTEL",   "The Indian Hotels Company Limited",                         "500850", 0.8100000000000001, 621.15, 1139040, "day", "21/08/2024 11:36:40"),
    @(9,  "PFC",        "Power Finance Corporation Limited",                         "532810", -1.06,              515.65, 7836483, "day", "21/08/2024 11:36:40"),
    @(10, "EXIDEIND",   "Exide Industries Limited",                                  "500086", 2.16,                  508, 4693036, "day", "21/08/2024 11:36:40"),
    @(11, "CROMPTON",   "Crompton Greaves Consumer Electricals Limited",             "539876", 3.19,                  468, 5377621, "day", "21/08/2024 11:36:40"),
    @(12, "HINDCOPPER", "Hindustan Copper Limited",                                  "513599", -0.52,                 317, 3062167, "day", "21/08/2024 11:36:40")
)

$startRow = 419
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # bsecode must stay TEXT even though it looks numeric - force the "@"
    # text format before assignment, then strip the format back off so the
    # cell keeps the workbook's default (unstyled) look, same as every
    # other data cell.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]
    $dCell.Style = "Normal"

    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
